$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 16.12468433333333
$ws.Range("H2").Value = 48.374053
$ws.Range("I2").Value = 0.2955490655206278
$ws.Range("J2").Value = 0.2955490655206279
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 5.397241999999999
$ws.Range("N2").Value = 16.191726
$ws.Range("O2").Value = 0.4331003391330405
$ws.Range("P2").Value = 0.4331003391330406
$ws.Range("Q2").Value = 87.02882352060864
$ws.Range("R2").Value = 783.2594116854779
$ws.Range("S2").Value = 0.1280024005074371
$ws.Range("T2").Value = 0.1280024005074372

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 16.12468433333333
$ws.Range("H3").Value = 48.374053
$ws.Range("I3").Value = 0.2955490655206278
$ws.Range("J3").Value = 0.2955490655206279
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 4.372979
$ws.Range("N3").Value = 13.118937
$ws.Range("O3").Value = 0.3509086099755513
$ws.Range("P3").Value = 0.3509086099755513
$ws.Range("Q3").Value = 70.51290597129565
$ws.Range("R3").Value = 634.6161537416609
$ws.Range("S3").Value = 0.1037107117614167
$ws.Range("T3").Value = 0.1037107117614167

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 16.12468433333333
$ws.Range("H4").Value = 48.374053
$ws.Range("I4").Value = 0.2955490655206278
$ws.Range("J4").Value = 0.2955490655206279
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.59901
$ws.Range("N4").Value = 1.79703
$ws.Range("O4").Value = 0.04806740815847847
$ws.Range("P4").Value = 0.04806740815847847
$ws.Range("Q4").Value = 9.65884716251
$ws.Range("R4").Value = 86.92962446259
$ws.Range("S4").Value = 0.01420627756323691
$ws.Range("T4").Value = 0.01420627756323692

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 16.12468433333333
$ws.Range("H5").Value = 48.374053
$ws.Range("I5").Value = 0.2955490655206278
$ws.Range("J5").Value = 0.2955490655206279
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.092643333333333
$ws.Range("N5").Value = 6.27793
$ws.Range("O5").Value = 0.1679236427329297
$ws.Range("P5").Value = 0.1679236427329297
$ws.Range("Q5").Value = 33.74321317225444
$ws.Range("R5").Value = 303.6889185502899
$ws.Range("S5").Value = 0.04962967568853715
$ws.Range("T5").Value = 0.04962967568853716

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 14.68975
$ws.Range("H6").Value = 44.06925
$ws.Range("I6").Value = 0.2692481784748309
$ws.Range("J6").Value = 0.2692481784748309
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 5.397241999999999
$ws.Range("N6").Value = 16.191726
$ws.Range("O6").Value = 0.4331003391330405
$ws.Range("P6").Value = 0.4331003391330406
$ws.Range("Q6").Value = 79.28413566949999
$ws.Range("R6").Value = 713.5572210254999
$ws.Range("S6").Value = 0.1166114774084027
$ws.Range("T6").Value = 0.1166114774084027

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 14.68975
$ws.Range("H7").Value = 44.06925
$ws.Range("I7").Value = 0.2692481784748309
$ws.Range("J7").Value = 0.2692481784748309
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.372979
$ws.Range("N7").Value = 13.118937
$ws.Range("O7").Value = 0.3509086099755513
$ws.Range("P7").Value = 0.3509086099755513
$ws.Range("Q7").Value = 64.23796826524999
$ws.Range("R7").Value = 578.1417143872499
$ws.Range("S7").Value = 0.09448150404705208
$ws.Range("T7").Value = 0.09448150404705206

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 14.68975
$ws.Range("H8").Value = 44.06925
$ws.Range("I8").Value = 0.2692481784748309
$ws.Range("J8").Value = 0.2692481784748309
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.59901
$ws.Range("N8").Value = 1.79703
$ws.Range("O8").Value = 0.04806740815847847
$ws.Range("P8").Value = 0.04806740815847847
$ws.Range("Q8").Value = 8.7993071475
$ws.Range("R8").Value = 79.19376432750001
$ws.Range("S8").Value = 0.01294206209067655
$ws.Range("T8").Value = 0.01294206209067656

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 14.68975
$ws.Range("H9").Value = 44.06925
$ws.Range("I9").Value = 0.2692481784748309
$ws.Range("J9").Value = 0.2692481784748309
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.092643333333333
$ws.Range("N9").Value = 6.27793
$ws.Range("O9").Value = 0.1679236427329297
$ws.Range("P9").Value = 0.1679236427329297
$ws.Range("Q9").Value = 30.74040740583333
$ws.Range("R9").Value = 276.6636666525
$ws.Range("S9").Value = 0.04521313492869961
$ws.Range("T9").Value = 0.04521313492869961

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.967860666666667
$ws.Range("H10").Value = 5.903582
$ws.Range("I10").Value = 0.03606888476606249
$ws.Range("J10").Value = 0.03606888476606249
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 5.397241999999999
$ws.Range("N10").Value = 16.191726
$ws.Range("O10").Value = 0.4331003391330405
$ws.Range("P10").Value = 0.4331003391330406
$ws.Range("Q10").Value = 10.62102024028133
$ws.Range("R10").Value = 95.589182162532
$ws.Range("S10").Value = 0.01562144622433222
$ws.Range("T10").Value = 0.01562144622433222

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.967860666666667
$ws.Range("H11").Value = 5.903582
$ws.Range("I11").Value = 0.03606888476606249
$ws.Range("J11").Value = 0.03606888476606249
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 4.372979
$ws.Range("N11").Value = 13.118937
$ws.Range("O11").Value = 0.3509086099755513
$ws.Range("P11").Value = 0.3509086099755513
$ws.Range("Q11").Value = 8.605413370259333
$ws.Range("R11").Value = 77.44872033233399
$ws.Range("S11").Value = 0.01265688221662733
$ws.Range("T11").Value = 0.01265688221662732

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.967860666666667
$ws.Range("H12").Value = 5.903582
$ws.Range("I12").Value = 0.03606888476606249
$ws.Range("J12").Value = 0.03606888476606249
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.59901
$ws.Range("N12").Value = 1.79703
$ws.Range("O12").Value = 0.04806740815847847
$ws.Range("P12").Value = 0.04806740815847847
$ws.Range("Q12").Value = 1.17876821794
$ws.Range("R12").Value = 10.60891396146
$ws.Range("S12").Value = 0.001733737805871452
$ws.Range("T12").Value = 0.001733737805871452

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.967860666666667
$ws.Range("H13").Value = 5.903582
$ws.Range("I13").Value = 0.03606888476606249
$ws.Range("J13").Value = 0.03606888476606249
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 2.092643333333333
$ws.Range("N13").Value = 6.27793
$ws.Range("O13").Value = 0.1679236427329297
$ws.Range("P13").Value = 0.1679236427329297
$ws.Range("Q13").Value = 4.118030505028889
$ws.Range("R13").Value = 37.06227454526
$ws.Range("S13").Value = 0.006056818519231489
$ws.Range("T13").Value = 0.006056818519231489

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 21.77610566666667
$ws.Range("H14").Value = 65.328317
$ws.Range("I14").Value = 0.3991338712384788
$ws.Range("J14").Value = 0.3991338712384788
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 5.397241999999999
$ws.Range("N14").Value = 16.191726
$ws.Range("O14").Value = 0.4331003391330405
$ws.Range("P14").Value = 0.4331003391330406
$ws.Range("Q14").Value = 117.5309121005713
$ws.Range("R14").Value = 1057.778208905142
$ws.Range("S14").Value = 0.1728650149928685
$ws.Range("T14").Value = 0.1728650149928685

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 21.77610566666667
$ws.Range("H15").Value = 65.328317
$ws.Range("I15").Value = 0.3991338712384788
$ws.Range("J15").Value = 0.3991338712384788
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 4.372979
$ws.Range("N15").Value = 13.118937
$ws.Range("O15").Value = 0.3509086099755513
$ws.Range("P15").Value = 0.3509086099755513
$ws.Range("Q15").Value = 95.22645278211434
$ws.Range("R15").Value = 857.038075039029
$ws.Range("S15").Value = 0.1400595119504553
$ws.Range("T15").Value = 0.1400595119504552

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 21.77610566666667
$ws.Range("H16").Value = 65.328317
$ws.Range("I16").Value = 0.3991338712384788
$ws.Range("J16").Value = 0.3991338712384788
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.59901
$ws.Range("N16").Value = 1.79703
$ws.Range("O16").Value = 0.04806740815847847
$ws.Range("P16").Value = 0.04806740815847847
$ws.Range("Q16").Value = 13.04410505539
$ws.Range("R16").Value = 117.39694549851
$ws.Range("S16").Value = 0.01918533069869355
$ws.Range("T16").Value = 0.01918533069869355

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 21.77610566666667
$ws.Range("H17").Value = 65.328317
$ws.Range("I17").Value = 0.3991338712384788
$ws.Range("J17").Value = 0.3991338712384788
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 2.092643333333333
$ws.Range("N17").Value = 6.27793
$ws.Range("O17").Value = 0.1679236427329297
$ws.Range("P17").Value = 0.1679236427329297
$ws.Range("Q17").Value = 45.56962234931222
$ws.Range("R17").Value = 410.12660114381
$ws.Range("S17").Value = 0.06702401359646148
$ws.Range("T17").Value = 0.06702401359646148

